$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1888888888888889
$ws.Range("C2").Value = 0.5888888888888889
$ws.Range("J2").Value = 0.02222222222222222
$ws.Range("P2").Value = 0.1259259259259259
$ws.Range("S2").Value = 0.07407407407407407
$ws.Range("B3").Value = 0.006097560975609756
$ws.Range("C3").Value = 0.02439024390243903
$ws.Range("J3").Value = 0.01829268292682927
$ws.Range("P3").Value = 0.7378048780487805
$ws.Range("S3").Value = 0.2134146341463415
$ws.Range("B6").Value = 0.0502092050209205
$ws.Range("D6").Value = 0.008368200836820083
$ws.Range("F6").Value = 0.06694560669456066
$ws.Range("J6").Value = 0.2343096234309623
$ws.Range("O6").Value = 0.008368200836820083
$ws.Range("Q6").Value = 0.1213389121338912
$ws.Range("R6").Value = 0.05857740585774059
$ws.Range("S6").Value = 0.4518828451882845
$ws.Range("B7").Value = 0.0861244019138756
$ws.Range("D7").Value = 0.03349282296650718
$ws.Range("F7").Value = 0.06698564593301436
$ws.Range("J7").Value = 0.138755980861244
$ws.Range("O7").Value = 0.009569377990430622
$ws.Range("Q7").Value = 0.1578947368421053
$ws.Range("R7").Value = 0.05263157894736842
$ws.Range("S7").Value = 0.4545454545454545
$ws.Range("B8").Value = 0.07983193277310924
$ws.Range("D8").Value = 0.01050420168067227
$ws.Range("E8").Value = 0.002100840336134454
$ws.Range("F8").Value = 0.07352941176470588
$ws.Range("J8").Value = 0.1029411764705882
$ws.Range("O8").Value = 0.02310924369747899
$ws.Range("Q8").Value = 0.1848739495798319
$ws.Range("R8").Value = 0.09873949579831932
$ws.Range("S8").Value = 0.4243697478991597
$ws.Range("B9").Value = 0.09793814432989691
$ws.Range("D9").Value = 0.02061855670103093
$ws.Range("F9").Value = 0.08247422680412371
$ws.Range("J9").Value = 0.08762886597938144
$ws.Range("O9").Value = 0.02061855670103093
$ws.Range("Q9").Value = 0.1958762886597938
$ws.Range("R9").Value = 0.07731958762886598
$ws.Range("S9").Value = 0.4175257731958763
$ws.Range("B10").Value = 0.09506263817243921
$ws.Range("D10").Value = 0.02210759027266028
$ws.Range("E10").Value = 0.001473839351510685
$ws.Range("F10").Value = 0.0707442888725129
$ws.Range("J10").Value = 0.08769344141488578
$ws.Range("O10").Value = 0.01694915254237288
$ws.Range("Q10").Value = 0.2137067059690494
$ws.Range("R10").Value = 0.08916728076639646
$ws.Range("S10").Value = 0.4030950626381725
$ws.Range("G11").Value = 0.1010452961672474
$ws.Range("J11").Value = 0.10801393728223
$ws.Range("K11").Value = 0.156794425087108
$ws.Range("L11").Value = 0.627177700348432
$ws.Range("S11").Value = 0.006968641114982578
$ws.Range("G12").Value = 0.7527472527472527
$ws.Range("J12").Value = 0.2197802197802198
$ws.Range("K12").Value = 0.01098901098901099
$ws.Range("L12").Value = 0.005494505494505495
$ws.Range("S12").Value = 0.01098901098901099
$ws.Range("G13").Value = 0.6911764705882353
$ws.Range("J13").Value = 0.2941176470588235
$ws.Range("S13").Value = 0.01470588235294118
$ws.Range("F15").Value = 0.01801801801801802
$ws.Range("H15").Value = 0.1396396396396396
$ws.Range("I15").Value = 0.06756756756756757
$ws.Range("J15").Value = 0.3873873873873874
$ws.Range("K15").Value = 0.04954954954954955
$ws.Range("M15").Value = 0.01801801801801802
$ws.Range("N15").Value = 0.004504504504504504
$ws.Range("O15").Value = 0.08558558558558559
$ws.Range("S15").Value = 0.2297297297297297
$ws.Range("F16").Value = 0.02339181286549707
$ws.Range("H16").Value = 0.1695906432748538
$ws.Range("I16").Value = 0.0935672514619883
$ws.Range("J16").Value = 0.4152046783625731
$ws.Range("K16").Value = 0.08187134502923976
$ws.Range("M16").Value = 0.03508771929824561
$ws.Range("O16").Value = 0.01754385964912281
$ws.Range("S16").Value = 0.1637426900584795
$ws.Range("F17").Value = 0.01680672268907563
$ws.Range("H17").Value = 0.2079831932773109
$ws.Range("I17").Value = 0.09453781512605042
$ws.Range("J17").Value = 0.4033613445378151
$ws.Range("K17").Value = 0.1029411764705882
$ws.Range("M17").Value = 0.03151260504201681
$ws.Range("O17").Value = 0.04411764705882353
$ws.Range("S17").Value = 0.09873949579831932
$ws.Range("F18").Value = 0.004807692307692308
$ws.Range("H18").Value = 0.2019230769230769
$ws.Range("I18").Value = 0.1298076923076923
$ws.Range("J18").Value = 0.4519230769230769
$ws.Range("K18").Value = 0.05288461538461538
$ws.Range("M18").Value = 0.03365384615384615
$ws.Range("O18").Value = 0.05288461538461538
$ws.Range("S18").Value = 0.07211538461538461
$ws.Range("F19").Value = 0.01670297748729121
$ws.Range("H19").Value = 0.2011619462599855
$ws.Range("I19").Value = 0.06753812636165578
$ws.Range("J19").Value = 0.4059549745824256
$ws.Range("K19").Value = 0.112563543936093
$ws.Range("M19").Value = 0.02687000726216413
$ws.Range("N19").Value = 0.001452432824981845
$ws.Range("O19").Value = 0.06899055918663761
$ws.Range("S19").Value = 0.09876543209876543
